$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product names (column A) and on-hand quantities (column B) with the
# refreshed finished-goods stock figures.
$updates = @(
    @{Row=3; A=$null; B=32213},
    @{Row=4; A='Сб. Фитонефрол (Урологический сбор) 50г'; B=5180},
    @{Row=5; A='Ромашка цветки вн 50г'; B=65730},
    @{Row=6; A='Береза почки 50г'; B=12611},
    @{Row=7; A='Ноготки цветки 50г'; B=17401},
    @{Row=8; A=$null; B=16786},
    @{Row=9; A='Эрва шерстистая трава 30г'; B=9895},
    @{Row=10; A='Багульник болотный побеги 50г'; B=11219},
    @{Row=11; A='Подорожник большой листья 50г'; B=8123},
    @{Row=12; A='Боярышник плоды 75г'; B=18950},
    @{Row=13; A='Пустырник трава 50г'; B=10346},
    @{Row=14; A='Рябина плоды 50г'; B=1554},
    @{Row=15; A='Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г'; B=3821},
    @{Row=16; A='Шалфей листья 50г'; B=33684},
    @{Row=17; A='Сб. Грудной №4 50г'; B=32466},
    @{Row=18; A='Бессмертник песчаный цветки 30г'; B=26575},
    @{Row=19; A='Девясил корневища и корни 50г'; B=16795},
    @{Row=20; A='Чабрец трава 50г'; B=19096},
    @{Row=21; A=$null; B=18968},
    @{Row=22; A='Аир корневища 75г'; B=7255},
    @{Row=23; A='Лен семена 100г'; B=57906},
    @{Row=24; A='Ламинарии слоевища (морская капуста) 100г'; B=16068},
    @{Row=25; A='Мать-и-мачеха листья 35г'; B=28340},
    @{Row=26; A='Укроп пахучий плоды 50г'; B=67012},
    @{Row=27; A='Крушина кора 50г'; B=11136},
    @{Row=28; A='Пижма цветки 75г'; B=18692},
    @{Row=29; A='Полынь горькая трава 50г'; B=47150},
    @{Row=30; A='Череда трава 50г'; B=13549},
    @{Row=31; A='Брусника листья 50г'; B=19079},
    @{Row=32; A='Шиповник плоды низковитаминные 50г'; B=40804},
    @{Row=33; A='Тысячелистник трава 50г'; B=17389},
    @{Row=34; A='Липа цветки 35г'; B=26228},
    @{Row=35; A='Зверобой трава 50г'; B=42098},
    @{Row=36; A='Эвкалипт прутовидный листья 75г'; B=35215},
    @{Row=37; A='Солодка корни 50г'; B=45935},
    @{Row=38; A='Кукуруза столбики с рыльцами 40г'; B=35049},
    @{Row=39; A='Можжевельник плоды 50г'; B=17400},
    @{Row=40; A='Сб. Фитопектол №1 (Грудной сбор №1) 35г'; B=7363},
    @{Row=41; A='Толокнянка листья 50г'; B=10754},
    @{Row=42; A='Сенна листья 50г'; B=32437},
    @{Row=43; A='Чага (березовый гриб) 50г'; B=41272},
    @{Row=44; A='Алтей корни 75г'; B=8759},
    @{Row=45; A=$null; B=21247},
    @{Row=46; A='Крапива листья 50г'; B=21593},
    @{Row=47; A='Чистотел трава 50г'; B=27622},
    @{Row=48; A=$null; B=10200},
    @{Row=49; A=$null; B=350},
    @{Row=50; A=$null; B=11676},
    @{Row=51; A='Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г'; B=3400},
    @{Row=52; A='Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г'; B=4380},
    @{Row=53; A=$null; B=4770},
    @{Row=54; A=$null; B=7830},
    @{Row=55; A=$null; B=9260},
    @{Row=60; A=$null; B=972},
    @{Row=61; A=$null; B=85247},
    @{Row=62; A='Фп Сб. Грудной №4 20x2,0г'; B=512496},
    @{Row=63; A='Фп Сенна листья 20x1,5г'; B=49400},
    @{Row=64; A='Фп Сб. Бруснивер 20x2,0г'; B=169879},
    @{Row=65; A='Фп Пижма цветки 20х1,5г'; B=4764},
    @{Row=66; A='Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г'; B=72211},
    @{Row=67; A='Фп Липа цветки 20x1,5г'; B=59181},
    @{Row=68; A='Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г'; B=74643},
    @{Row=69; A='Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г'; B=1692},
    @{Row=70; A='Фп Фиточай "Опалиховский" (БАД) 20х2,0 г'; B=4680},
    @{Row=71; A='Фп Фиточай "Тибетский" (БАД) 20х2,0  г'; B=9683},
    @{Row=72; A='Фп Мята перечная листья 20x1,5г'; B=63928},
    @{Row=73; A='Фп Сб. Арфазетин-Э 20x2,0г'; B=38789},
    @{Row=74; A='Фп Сб. Элекасол 20x2,0г'; B=39834},
    @{Row=75; A='Фп Чистотел трава 20х1,5г'; B=30732},
    @{Row=76; A='Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г'; B=1746},
    @{Row=77; A='Фп Брусника листья 20х1,5г'; B=77597},
    @{Row=78; A='Фп Подорожник листья 20x1,5г'; B=30230},
    @{Row=79; A='Фп Ромашка цветки 20x1,5г'; B=1459368},
    @{Row=80; A='Фп Пустырник трава 20x1,5г'; B=43950},
    @{Row=81; A='Фп Пастушья сумка трава 20х1,5г'; B=5398},
    @{Row=82; A='Фп Шиповник плоды 20х2,0г'; B=53820},
    @{Row=83; A='Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г'; B=25162},
    @{Row=84; A='Фп Череда трава 20х1,5г'; B=52793},
    @{Row=85; A='Фп Мелисса лекарственная трава 20x1,5г'; B=42174},
    @{Row=86; A='Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г'; B=204529},
    @{Row=87; A='Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г'; B=96255},
    @{Row=88; A='Фп Зверобой трава 20x1,5г'; B=59657},
    @{Row=89; A='Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г'; B=63667},
    @{Row=90; A='Фп Чабрец трава 20x1,5 г'; B=81702},
    @{Row=91; A='Фп Крапива листья 20x1,5г'; B=78545},
    @{Row=92; A='Фп Хвощ полевой трава 20х1,5г'; B=33658},
    @{Row=93; A='Фп Душица трава 20x1,5г'; B=34722},
    @{Row=94; A='Фп Сб. Желудочный №3 20x2,0г'; B=29191},
    @{Row=95; A='Фп Толокнянка листья 20x1,5г'; B=48562},
    @{Row=96; A='Фп Береза листья 20x1,5г'; B=6160},
    @{Row=97; A='Фп Золототысячник трава 20х1,5г'; B=5997},
    @{Row=98; A='Фп Боярышник плоды 20х3,0г'; B=26576},
    @{Row=99; A='Фп Фиалка трехцветная трава 20x1,5г'; B=5656},
    @{Row=100; A='Фп Аир корневища 20x1,5г'; B=6863},
    @{Row=101; A=$null; B=6001},
    @{Row=103; A='Фп Лапчатка корневища 20x2,5г'; B=2659},
    @{Row=104; A='Фп Тысячелистник трава 20x1,5г'; B=13754},
    @{Row=105; A='Фп Крушина кора 20x1,5г'; B=7533},
    @{Row=106; A='Фп Ноготки цветки 20x1,5г'; B=34763},
    @{Row=107; A=$null; B=12870},
    @{Row=108; A='Фп Бадан корневища 20x1,5г'; B=1807},
    @{Row=109; A='Фп Валериана корневища с корнями 20x1,5г'; B=34262},
    @{Row=110; A=$null; B=11752},
    @{Row=111; A=$null; B=166941}
)

foreach ($u in $updates) {
    if ($null -ne $u.A) {
        $ws.Cells.Item($u.Row, 1).Value = $u.A
    }
    if ($null -ne $u.B) {
        $ws.Cells.Item($u.Row, 2).Value = $u.B
    }
}

# Restore the saved view state (scroll position / active selection).
$ws.Activate()
$ws.Range("A92").Select()
$excel.ActiveWindow.ScrollRow = 76
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "Applied finish goods stock refresh."
